# "Add files via upload" -- refresh the ELM forecasting results on the
# DKI2 sheet: column B / column C values are recomputed, the stray extra
# cell style used only by B2 is retired (B2 now shares the plain bordered
# style already used throughout the table), and the saved selection /
# window state reflect where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DKI2")

# New C2:C32 values (column C is being recomputed).
$newC = @{
    2  = 31.246926210000002
    3  = 32.464636249999998
    4  = 35.1441199
    5  = 36.81661398
    6  = 32.838756609999997
    7  = 31.73053827
    8  = 33.694716309999997
    9  = 33.779772960000003
    10 = 32.705885979999998
    11 = 31.984211080000001
    12 = 31.181223630000002
    13 = 31.288028440000001
    14 = 32.233465639999999
    15 = 32.1091689
    16 = 31.529856980000002
    17 = 31.752060369999999
    18 = 33.437345460000003
    19 = 33.953034770000002
    20 = 32.989783359999997
    21 = 32.528868750000001
    22 = 32.867997090000003
    23 = 32.840477640000003
    24 = 32.607683979999997
    25 = 32.173532999999999
    26 = 31.729521770000002
    27 = 31.861761220000002
    28 = 32.369893779999998
    29 = 32.556455249999999
    30 = 32.15151007
    31 = 32.228290020000003
    32 = 32.79541184
}

# B2 picks up the value that used to sit in C2, and loses its one-off
# number-format style in favor of the plain bordered style already used
# by every other cell in the column (same look, one fewer cellXf record).
$ws.Range("B2").Value = 45.007682665611007
$ws.Range("B2").Borders.LineStyle = 1

foreach ($row in 2..32) {
    $ws.Cells.Item($row, 3).Value = $newC[$row]
}

# Author's last selection before saving.
[void]$ws.Range("H10").Select()
